# The sheet currently has an extra leading column A (a row-index style
# column) that is not part of the final data. Deleting it shifts every
# other column (B:F -> A:E) one position to the left, which matches the
# target layout (dimension A1:E3 instead of A1:F3).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("A").Delete()

# After the shift, the header that used to read "MODEL_CONDITION" (now in
# column D) needs its text corrected to "MODELCONDITION".
$ws.Cells.Item(1, 4).Value = "MODELCONDITION"
